$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44600
$ws.Range("J2").Value = 1300
$ws.Range("K2").Value = 3500
$ws.Range("L2").Value = 4000
$ws.Range("M2").Value = 3808
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 38
$ws.Range("D3").Value = 44232
$ws.Range("J3").Value = 16000
$ws.Range("D4").Value = 44168
$ws.Range("J4").Value = 7000
$ws.Range("D5").Value = 44167
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 7000
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = 3000
$ws.Range("P5").Value = 30
$ws.Range("D6").Value = 44229
$ws.Range("J6").Value = 16000
$ws.Range("D7").Value = 44188
$ws.Range("D8").Value = 44186
$ws.Range("J8").Value = 10000
$ws.Range("D9").Value = 44231
$ws.Range("J9").Value = 12000
$ws.Range("D10").Value = 44602
$ws.Range("J10").Value = 12000
$ws.Range("D11").Value = 44602
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 6000
$ws.Range("L11").Value = 2500
$ws.Range("M11").Value = 2500
$ws.Range("P11").Value = 25
$ws.Range("D12").Value = 44159
$ws.Range("D13").Value = 44230
$ws.Range("D14").Value = 44161
$ws.Range("J14").Value = 7000
$ws.Range("D15").Value = 44181
$ws.Range("J15").Value = 12000
$ws.Range("D16").Value = 44166
$ws.Range("J16").Value = 7000
$ws.Range("K16").Value = 3000
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = 3000
$ws.Range("O16").Value = "Provincia de Chacabuco"
$ws.Range("P16").Value = 30
$ws.Range("D17").Value = 44215
$ws.Range("J17").Value = 16000
$ws.Range("D18").Value = 44204
$ws.Range("D19").Value = 44245
$ws.Range("J19").Value = 9000
$ws.Range("K19").Value = 3000
$ws.Range("M19").Value = 3000
$ws.Range("O19").Value = "Región Metropolitana"
$ws.Range("P19").Value = 30
$ws.Range("D20").Value = 44245
$ws.Range("I20").Value = "Segunda"
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 2500
$ws.Range("L20").Value = 2500
$ws.Range("M20").Value = 2500
$ws.Range("O20").Value = "Región Metropolitana"
$ws.Range("P20").Value = 25
$ws.Range("D21").Value = 44160
$ws.Range("J21").Value = 7000
$ws.Range("D22").Value = 44162
$ws.Range("D23").Value = 44209
$ws.Range("K23").Value = 2500
$ws.Range("M23").Value = 2750
$ws.Range("P23").Value = 28
$ws.Range("D24").Value = 44214
$ws.Range("J24").Value = 7000
$ws.Range("O24").Value = "Provincia de Chacabuco"
$ws.Range("D25").Value = 44189
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 16000
$ws.Range("K25").Value = 3000
$ws.Range("L25").Value = 3000
$ws.Range("M25").Value = 3000
$ws.Range("O25").Value = "Provincia de Chacabuco"
$ws.Range("P25").Value = 30
$ws.Range("D26").Value = 44210
$ws.Range("J26").Value = 8800
$ws.Range("K26").Value = 2500
$ws.Range("M26").Value = 2750
$ws.Range("P26").Value = 28
$ws.Range("D27").Value = 44187
$ws.Range("J27").Value = 12000
